# Update column G ("K") values for rows 2-18 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 7
    3  = 2
    4  = 5
    5  = 2
    6  = 1
    7  = 8
    8  = 5
    9  = 1
    10 = 5
    11 = 2
    12 = 8
    13 = 7
    14 = 3
    15 = 4
    16 = 4
    17 = 0
    18 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
